# Add a new contact record (row 5) to the NewArrivals contact list, mirroring
# the existing rows 2-4: firstName, lastName, email (+hyperlink), address1,
# city, state (+hyperlink display text), zip.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row values -------------------------------------------------------
$ws.Range("A5").Value = "neoguest2"
$ws.Range("B5").Value = "abc"
$ws.Range("C5").Value = "neoguest2ab@test.com"
$ws.Range("E5").Value = "PLANT CITY"
$ws.Range("D5").Value = "1908 INDUSTRIAL PARK DR"
$ws.Range("F5").Value = "FLORIDA"
$ws.Range("G5").Value = 33566

# --- hyperlinks on the email and state cells (matches rows 2-4 pattern) ---
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:neoguest2ab@test.com") | Out-Null

$stateLink = $ws.Hyperlinks.Add($ws.Range("F5"), "mailto:C@bi`$ush5")
$stateLink.TextToDisplay = "C@bi`$ush5"
# Re-assert the real cell text: Excel uses TextToDisplay to (re)write the
# cell's text, but the source data keeps "FLORIDA" in the cell while only the
# hyperlink's stored display string is "C@bi$ush5".
$ws.Range("F5").Value = "FLORIDA"

# --- formatting: reuse the same styles as the rows above -------------------
$ws.Range("A2:C2").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("F2").Copy() | Out-Null
$ws.Range("F5").PasteSpecial(-4122) | Out-Null

# --- selection cursor, matching the saved view state ------------------------
$ws.Range("D9").Select() | Out-Null
